# Localization status update: file d2525536-afb4-49a1-9951-26382a7aad42.md
# moves from "Ready for handoff" back to "In Translation" for both
# locales (zh-cn, de-de) while a new handoff (HO Xliff) report is
# regenerated for archive.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for d2525536-... (row 4) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn detail sheet: Status column (C) for the same file, row 4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de detail sheet: Status column (C) for the same file, row 4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = "In Translation"
